$wb = $excel.ActiveWorkbook

# --- YDS sheet: append new week play-by-play yardage logs ---
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = "7 8 4 1 7 2 3 1 5 8 4 3 2 2 9 0 14 3 -6 7 15 0 18 3 0 3 5 3 2 6 1 1 0 2 5 -1 -1 16 6 0 7 3 5 9 0 2 -2 0 11 0 -1 1 3 2 -5 3 9 0 3 3 5 -3 10 4 9 1 2 1 0 4 13 7 2 2 -1 4 4 12 0 2 14 0 1 6 6 4 14 12 2 -1 1 -1 14 7 4 11 -2 2 6 7 3 3 0 7 -2 7 7 5 1 5 6 0 26 3 1 5 0 3 3 -4 1 2 9 -5 8 2 3 11 1 8 6 11 3 5 -1 2 5 3 2 8 5 3 6 2 3 5 -2 6 9 6 0 5 2 5 2 6 30 0 13 3 6 13 7 -3 2 -2 2 2 4 7 -4 0 1 3 -2 2 4 9 3 9 2 2 1 -2 6 2 5 2 7 2 13 6 6 3 3 6 20 9 2 2 6 4 6 6 6 3 6 6 -2 1 11 8 2 6 -1 1 5 3 2 1 2 4 2 0 4 4 4 -2 3 10 2 1 3 9 6 0 3 0 4 10 1 4 4 1 -1 2 11 11 12 10 6 2 4 12 1 16 15 10 2 5 4 1 1 4 0 5 2 10 4 -9 3 15 4 8 6 -4 -3 1 7 -2 8 6 3 0 -3 4 3 8 3 4 9 6 24 2 4 0 3 6 10 4 4 4 2 1 15 5 2 2 3 1 1 -3 0 5 -1 5 3 -1 11 8 -1 12 8 4 2 4 1 4 1 4 -5 -1 5 11 17 4 0 1 2 -2 2 4 1 11 4 21 3 -2 0 7 8 -4 5 2 2 0 2 30 3 4 1 13 -1 1 7 1 7 8 5 7 5 33 2 -2 1 6 2 2 -2 0 4 5 0 3 3 4 11"
$ydsWs.Range("C2").Value = "2 1 1 2 6 2 2 0 2 2 1 1 1 9 2 12 2 2 0 0 4 2 12 0 1 7 1 5 9 -7 8 6 5 6 6 4 2 -1 1 -2 7 1 4 3 1 1 3 1 0 2 2 19 0 8 7 5 7 1 4 0 1 2 6 1 11 5 5 47 14 0 14 0 0 10 5 2 18 8 1 5 6 4 1 3 0 1 4 4 1 4 14 3 6 1 6 7 18 -3 12 1 -3 2 -3 3 4 14 5 9 1 3 4 3 5 1 8 9 6 0 14 2 0 30 16 -1 -1 23 4 0 3 3 -1 -3 16 5 4 16 14 -1 7 2 2 0 1 9 3 0 -2 -2 7 11 3 4 2 3 6 2 2 -2 4 19 0 4 9 2 4 1 -1 9 0 5 2 18 2 8 -2 3 3 14 5 0 5 3 1 2 2 0 16 12 1 5 2 -1 6 2 13 3 3 2 8 4 11 12 2 -2 3 6 5 1 4 11 2 5 1 7 -6 1 2 9 4 0 -1 1 6 2 4 0 4 0 0 4 4 6 11 11 1 9 4 3 10 2 -2 4 -3 15 5 2 11 1 8 3 4 0 1 4 5 11 20 2 0 11 9 4 9 4 2 8 -2 0 14 3 11 0 0 7 0 3 2 6 3 1 5 0 2 2 3 2 5 -3 3 0 2 16 0 2 4 4 3 3 8 7 2 0 8 2 1 9 2 1 1 5 -1 0 2 6 5 4 11 0 4 -1 4 11 4 0 11 1 4 3 11 8 5 9 0 2 1 -2 1 11 11 4 2 3 5 3 8 2 16 6 3 0 2 2 26 3 4 1 -2 12 3 1 2 4 1 8 7 4 5 2 1 5 3 4 5 55 0 4 4 4 -2 9 1 2 3 -1 3 1 1 2 0 11 14 7 2 -3 16 -1 12 -1 0 -1"
$ydsWs.Range("B3").Value = "11 13 -5 8 1 9 26 5 27 15 11 57 14 4 5 11 22 17 1 7 12 7 3 -2 32 20 20 19 8 5 13 18 19 5 16 7 17 2 10 6 7 13 8 1 4 8 11 11 16 9 7 1 29 12 26 1 13 8 24 12 7 8 10 10 16 25 10 8 30 5 8 11 23 9 29 9 3 13 6 5 17 1 8 5 9 39 9 15 2 11 4 7 13 8 6 5 5 8 6 13 14 5 6 12 11 -1 8 5 7 0 7 5 24 2 8 9 4 7 25 9 0 16 18 17 5 15 14 1 4 4 21 6 4 7 4 2 9 1 3 5 5 12 21 5 12 6 13 7 5 1 7 7 19 4 17 -2 13 10 23 7 9 21 4 11 1 11 13 6 2 6 5 33 13 -3 6 21 19 17 6 9 7 12 13 17 2 6 16 -2 6 10 2 3 11 6 9 2 10 6 10 1 2 3 9 8 15 7 3 8 13 10 5 3 11 2 19 7 9 4 10 16 27 7 6 4 8 11 7 64 4 13 4 7 18 24 13 17 7 11 15 22 10 23 9 8 14 6 18 7 9 3 19 19 13 18 11 1 9 8 8 9 4 11 13 9 16 4 23 5 6 8 5 8 14 8 10 5 15 63 0 4 6 2 9 19 8 6 13 2 18 14 6 6 17 6"
$ydsWs.Range("C3").Value = "7 35 10 14 10 8 -3 1 11 20 7 8 22 6 15 2 40 12 8 9 7 4 16 3 7 4 27 3 8 23 -1 11 4 5 8 30 10 11 1 10 5 5 13 15 11 5 12 5 18 12 10 14 18 4 7 20 35 6 7 6 23 8 -7 6 -5 3 10 8 8 2 7 2 0 1 15 17 7 11 8 13 5 53 14 20 11 9 6 7 7 3 3 20 7 24 13 2 -6 8 14 7 9 24 35 5 17 26 40 3 4 12 6 8 6 7 13 -3 27 15 -2 6 18 7 1 27 15 7 16 15 1 6 8 15 -4 2 14 16 2 5 8 9 12 4 12 3 7 7 15 11 2 2 18 1 -3 2 9 17 5 8 15 10 7 8 12 -1 4 9 41 3 28 3 7 18 1 4 19 5 20 4 0 5 19 21 7 9 -5 17 5 3 30 7 39 6 4 5 12 12 7 13 4 15 10 18 2 5 -1 25 4 10 6 0 9 57 9 5 5 5 7 11 11 7 7 4 11 11 9 2 -2 5 6 -2 6 36 9 13 3 20 4 4 14 12 7 14 3 11 6 2 23 12 9 15 9 10 12 11 4 10 9 7 11 7 7 7 8 28 20 4 3 4 13 7 19 7 1 10 22 62 13 23 11 12 7 10"

# --- OFF sheet: updated season totals after Week 16 ---
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("C2").Value = 198
$offWs.Range("F2").Value = 47
$offWs.Range("G2").Value = 59
$offWs.Range("J2").Value = 35
$offWs.Range("L2").Value = 290
$offWs.Range("M2").Value = 167
$offWs.Range("O2").Value = 28
$offWs.Range("Q2").Value = 508
$offWs.Range("C3").Value = 163
$offWs.Range("E3").Value = 34
$offWs.Range("F3").Value = 121
$offWs.Range("H3").Value = 25
$offWs.Range("I3").Value = 71
$offWs.Range("J3").Value = 43
$offWs.Range("N3").Value = 20

# --- DEF sheet: updated season totals after Week 16 ---
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 197
$defWs.Range("D2").Value = 7
$defWs.Range("F2").Value = 63
$defWs.Range("G2").Value = 47
$defWs.Range("J2").Value = 30
$defWs.Range("L2").Value = 241
$defWs.Range("M2").Value = 150
$defWs.Range("Q2").Value = 499
$defWs.Range("B3").Value = 10
$defWs.Range("C3").Value = 128
$defWs.Range("D3").Value = 7
$defWs.Range("E3").Value = 35
$defWs.Range("F3").Value = 91
$defWs.Range("G3").Value = 22
$defWs.Range("H3").Value = 37
$defWs.Range("I3").Value = 45
$defWs.Range("J3").Value = 53

# --- ST sheet: special teams totals + appended per-game logs ---
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 67
$stWs.Range("D2").Value = 65
$stWs.Range("J2").Value = 55
$stWs.Range("K2").Value = 53
$stWs.Range("B3").Value = 43
$stWs.Range("D3").Value = "36 18 54 35 46 36 33 39 51 37 32 42 37 49 49 43 44 42 44 30 40 33 40 40 44 39 44 36 28 39 46 32 40 34 39 49 36 46 44 57 59 46 43 44 55 54 48 50 32 55 57 55 40 51 54 33 47 54 41 40 40 50 55 48 60"
$stWs.Range("D4").Value = "0 0 15 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 22 0 0 7 0 0 0 0 0 0 2 12 6 0 4 8 0 13 0 18 15 7 0 1 17 7 0 0 0 0 7 7 0 0 0 0 17 14 0 0 4 9 20 6 5"
$stWs.Range("D5").Value = "4 9 0 0 0 12 15 0 10 14 0 0 0 5 0 12 0 0 0 0 0 0 0 0 0 22 8 0 9 1 0 4 0 0 10 0 0 4 0 0 0 0 0 0 4 0 10 0 16 10 0 0 0 0 0 0 0 0 0 0 0 0 16 16 6 0"
$stWs.Range("B6").Value = "21 21 39 13 31 24 16 13 10 21 22 26 21 26 28 0 25 21 10 18 21 27 25 26"

# --- TURNS sheet: updated turnover totals ---
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("B2").Value = 11
$turnsWs.Range("D2").Value = 8
$turnsWs.Range("D3").Value = 7
$turnsWs.Range("E3").Value = 11

# --- PEN sheet: updated penalty totals ---
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("D2").Value = 12
$penWs.Range("D3").Value = 1

